$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 352.53845
$ws.Range("I28").Value = 363.36365
$ws.Range("J28").Value = 293
$ws.Range("K28").Value = 363.36365
$ws.Range("L28").Value = 293
$ws.Range("M28").Value = 121.63635
$ws.Range("N28").Value = -1263

# Row 40
$ws.Range("H40").Value = 13705.556
$ws.Range("I40").Value = 2510
$ws.Range("J40").Value = 27700
$ws.Range("K40").Value = 2510
$ws.Range("L40").Value = 27700
$ws.Range("M40").Value = -2335
$ws.Range("N40").Value = -28050

# Row 109
$ws.Range("H109").Value = 35992
$ws.Range("J109").Value = 35992
$ws.Range("L109").Value = 35992
$ws.Range("N109").Value = -38766

# Row 114
$ws.Range("H114").Value = 39372.668
$ws.Range("J114").Value = 39372.668
$ws.Range("L114").Value = 39372.668
$ws.Range("N114").Value = -48050.668

# Row 117
$ws.Range("H117").Value = 48403.332
$ws.Range("J117").Value = 48403.332
$ws.Range("L117").Value = 48403.332
$ws.Range("N117").Value = -57581.332

# Row 120
$ws.Range("H120").Value = 49722
$ws.Range("J120").Value = 49722
$ws.Range("L120").Value = 49722
$ws.Range("N120").Value = -59398

# Row 125
$ws.Range("H125").Value = 600.25
$ws.Range("I125").Value = 461.6
$ws.Range("J125").Value = 831.3333
$ws.Range("K125").Value = 4154.400000000001
$ws.Range("L125").Value = 7481.9997
$ws.Range("M125").Value = -1694.400000000001
$ws.Range("N125").Value = -12401.9997

# Row 132
$ws.Range("H132").Value = 13642.176
$ws.Range("I132").Value = 2057.9321
$ws.Range("J132").Value = 59206.867
$ws.Range("K132").Value = 6173.7963
$ws.Range("L132").Value = 177620.601
$ws.Range("M132").Value = -3643.7963
$ws.Range("N132").Value = -182680.601

# Row 138
$ws.Range("H138").Value = 1408.02
$ws.Range("I138").Value = 646.5
$ws.Range("J138").Value = 1704.1666
$ws.Range("K138").Value = 1939.5
$ws.Range("L138").Value = 5112.4998
$ws.Range("M138").Value = 3200.5
$ws.Range("N138").Value = -15392.4998

$ws = $wb.Worksheets.Item("ARM")
# Row 106
$ws.Range("H106").Value = 46840.5
$ws.Range("J106").Value = 46840.5
$ws.Range("L106").Value = 46840.5
$ws.Range("N106").Value = -49364.5

# Row 107
$ws.Range("H107").Value = 36970.75
$ws.Range("J107").Value = 36970.75
$ws.Range("L107").Value = 36970.75
$ws.Range("N107").Value = -44650.75

# Row 109
$ws.Range("H109").Value = 45091.75
$ws.Range("J109").Value = 45091.75
$ws.Range("L109").Value = 45091.75
$ws.Range("N109").Value = -47865.75

# Row 111
$ws.Range("H111").Value = 47495.5
$ws.Range("J111").Value = 47495.5
$ws.Range("L111").Value = 47495.5
$ws.Range("N111").Value = -55675.5

# Row 117
$ws.Range("H117").Value = 48408.6
$ws.Range("J117").Value = 48408.6
$ws.Range("L117").Value = 48408.6
$ws.Range("N117").Value = -57586.6

# Row 118
$ws.Range("H118").Value = 49621.5
$ws.Range("J118").Value = 49621.5
$ws.Range("L118").Value = 49621.5
$ws.Range("N118").Value = -52935.5

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 4000
$ws.Range("I22").Value = 4000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 4000
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -3827

# Row 110
$ws.Range("H110").Value = 48694
$ws.Range("J110").Value = 48694
$ws.Range("L110").Value = 48694
$ws.Range("N110").Value = -56874

# Row 119
$ws.Range("H119").Value = 48376.5
$ws.Range("J119").Value = 48376.5
$ws.Range("L119").Value = 48376.5
$ws.Range("N119").Value = -58052.5

# Row 120
$ws.Range("H120").Value = 48761
$ws.Range("J120").Value = 48761
$ws.Range("L120").Value = 48761
$ws.Range("N120").Value = -58437

$ws = $wb.Worksheets.Item("CRP")
# Row 110
$ws.Range("H110").Value = 43037.668
$ws.Range("J110").Value = 43037.668
$ws.Range("L110").Value = 43037.668
$ws.Range("N110").Value = -51217.668

# Row 112
$ws.Range("H112").Value = 26286.143
$ws.Range("J112").Value = 26286.143
$ws.Range("L112").Value = 26286.143
$ws.Range("N112").Value = -29240.143

# Row 116
$ws.Range("H116").Value = 47887
$ws.Range("J116").Value = 47887
$ws.Range("L116").Value = 47887
$ws.Range("N116").Value = -57065

$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Row 107
$ws.Range("H107").Value = 12255.059
$ws.Range("J107").Value = 17233.834
$ws.Range("L107").Value = 51701.50199999999
$ws.Range("N107").Value = -55541.50199999999

# Row 132
$ws.Range("H132").Value = 975
$ws.Range("I132").Value = 750
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 6750
$ws.Range("L132").Value = 10800
$ws.Range("M132").Value = -4220
$ws.Range("N132").Value = -15860

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 15416.667
$ws.Range("J5").Value = 15454.546
$ws.Range("L5").Value = 15454.546
$ws.Range("N5").Value = -15678.546

# Row 52
$ws.Range("H52").Value = 21500
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

# Row 110
$ws.Range("H110").Value = 46888
$ws.Range("J110").Value = 46888
$ws.Range("L110").Value = 46888
$ws.Range("N110").Value = -55068

# Row 122
$ws.Range("H122").Value = 1163.7646
$ws.Range("I122").Value = 1214.1538
$ws.Range("K122").Value = 3642.4614
$ws.Range("M122").Value = -1192.4614

# Row 123
$ws.Range("H123").Value = 13900
$ws.Range("J123").Value = 13900
$ws.Range("L123").Value = 13900
$ws.Range("N123").Value = -18800

# Row 132
$ws.Range("H132").Value = 3687.122
$ws.Range("I132").Value = 1583.2174
$ws.Range("J132").Value = 6375.4443
$ws.Range("K132").Value = 4749.6522
$ws.Range("L132").Value = 19126.3329
$ws.Range("M132").Value = -2219.6522
$ws.Range("N132").Value = -24186.3329

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 26222.555
$ws.Range("J2").Value = 46000.6
$ws.Range("L2").Value = 46000.6
$ws.Range("N2").Value = -46224.6

# Row 40
$ws.Range("H40").Value = 3132.6667
$ws.Range("I40").Value = 2460.5
$ws.Range("K40").Value = 2460.5
$ws.Range("M40").Value = -2324.5

# Row 110
$ws.Range("H110").Value = 45643
$ws.Range("J110").Value = 45643
$ws.Range("L110").Value = 45643
$ws.Range("N110").Value = -53823

# Row 120
$ws.Range("H120").Value = 52090.5
$ws.Range("J120").Value = 52090.5
$ws.Range("L120").Value = 52090.5
$ws.Range("N120").Value = -61766.5

# Row 123
$ws.Range("H123").Value = 39421
$ws.Range("J123").Value = 39421
$ws.Range("L123").Value = 39421
$ws.Range("N123").Value = -49221

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 800
$ws.Range("I2").Value = 766.6667
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 766.6667
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -654.6667
$ws.Range("N2").Value = -1224

# Row 108
$ws.Range("H108").Value = 28446
$ws.Range("J108").Value = 28446
$ws.Range("L108").Value = 28446
$ws.Range("N108").Value = -36126

# Row 110
$ws.Range("H110").Value = 23386.4
$ws.Range("J110").Value = 23386.4
$ws.Range("L110").Value = 23386.4
$ws.Range("N110").Value = -31566.4
